$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 10 updates
$ws.Range("C10").Value = 13
$ws.Range("F10").Value = "Haris"

# Row 25 updates - new task row
$ws.Range("A25").Value = "Make a simple Adobe Xd design"
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "Haris"
$ws.Range("F25").Value = "Haris"

# Update selection to F10 as in diff
$ws.Range("F10").Select()
